$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2554.75
$ws.Range("J17").Value = 2554.75
$ws.Range("L17").Value = 7664.25
$ws.Range("N17").Value = -8000.25

$ws.Range("H28").Value = 2128.2856
$ws.Range("I28").Value = 479.6
$ws.Range("K28").Value = 479.6
$ws.Range("M28").Value = 5.399999999999977

$ws.Range("H64").Value = 4499
$ws.Range("I64").Value = 4998
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 4998
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -4750
$ws.Range("N64").Value = -4496

$ws.Range("H67").Value = 4499
$ws.Range("I67").Value = 4998
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 4998
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -4140
$ws.Range("N67").Value = -5716

$ws.Range("H92").Value = 1513.3334
$ws.Range("I92").Value = 1647.1428
$ws.Range("K92").Value = 1647.1428
$ws.Range("M92").Value = -399.1428000000001

$ws.Range("H112").Value = 1820.5625
$ws.Range("I112").Value = 1324.75
$ws.Range("J112").Value = 1985.8334
$ws.Range("K112").Value = 3974.25
$ws.Range("L112").Value = 5957.5002
$ws.Range("M112").Value = -2866.25
$ws.Range("N112").Value = -8173.5002

$ws.Range("H125").Value = 3766.125
$ws.Range("J125").Value = 8466.333000000001
$ws.Range("L125").Value = 76196.997
$ws.Range("N125").Value = -81116.997

$ws.Range("H131").Value = 500
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 500
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 1500
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -11580

$ws.Range("H132").Value = 2297.762
$ws.Range("I132").Value = 2297.762
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6893.286
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4363.286
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7037.647
$ws.Range("I32").Value = 7008.75
$ws.Range("J32").Value = 7500
$ws.Range("K32").Value = 7008.75
$ws.Range("L32").Value = 7500
$ws.Range("M32").Value = -6721.75
$ws.Range("N32").Value = -8074

$ws.Range("H63").Value = 2991.6667
$ws.Range("I63").Value = 2487.5
$ws.Range("K63").Value = 2487.5
$ws.Range("M63").Value = -1801.5

$ws.Range("H66").Value = 2991.6667
$ws.Range("I66").Value = 2487.5
$ws.Range("K66").Value = 12437.5
$ws.Range("M66").Value = -9005.5

$ws.Range("H88").Value = 3088.75
$ws.Range("J88").Value = 4585.4287
$ws.Range("L88").Value = 4585.4287
$ws.Range("N88").Value = -5397.4287

$ws.Range("H91").Value = 3088.75
$ws.Range("J91").Value = 4585.4287
$ws.Range("L91").Value = 4585.4287
$ws.Range("N91").Value = -7393.4287

$ws.Range("H132").Value = 3337.3333
$ws.Range("I132").Value = 2012
$ws.Range("K132").Value = 6036
$ws.Range("M132").Value = -3506

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2347.1428
$ws.Range("I86").Value = 1757.5
$ws.Range("J86").Value = 3133.3333
$ws.Range("K86").Value = 1757.5
$ws.Range("L86").Value = 3133.3333
$ws.Range("M86").Value = -634.5
$ws.Range("N86").Value = -5379.3333

$ws.Range("H89").Value = 2347.1428
$ws.Range("I89").Value = 1757.5
$ws.Range("J89").Value = 3133.3333
$ws.Range("K89").Value = 8787.5
$ws.Range("L89").Value = 15666.6665
$ws.Range("M89").Value = -3171.5
$ws.Range("N89").Value = -26898.6665

$ws.Range("H105").Value = 2500
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws.Range("H132").Value = 250750
$ws.Range("J132").Value = 250750
$ws.Range("L132").Value = 250750
$ws.Range("N132").Value = -260870

$ws.Range("H134").Value = 4186.3076
$ws.Range("I134").Value = 4186.3076
$ws.Range("K134").Value = 12558.9228
$ws.Range("M134").Value = -10023.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 27027.666
$ws.Range("I50").Value = 21083
$ws.Range("K50").Value = 21083
$ws.Range("M50").Value = -20458

$ws.Range("H62").Value = 4799.8
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 5999.5
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 5999.5
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -7247.5

$ws.Range("H65").Value = 4799.8
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 5999.5
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 29997.5
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -36237.5

$ws.Range("H132").Value = 1217.7
$ws.Range("I132").Value = 1217.7
$ws.Range("K132").Value = 3653.1
$ws.Range("M132").Value = -1123.1

$ws.Range("H134").Value = 2874.261
$ws.Range("I134").Value = 2766.8572
$ws.Range("K134").Value = 8300.571599999999
$ws.Range("M134").Value = -5765.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 133.33333
$ws.Range("J33").Value = 200
$ws.Range("L33").Value = 1200
$ws.Range("N33").Value = -1766

$ws.Range("H68").Value = 1584.2858
$ws.Range("J68").Value = 1532.6666
$ws.Range("L68").Value = 4597.9998
$ws.Range("N68").Value = -6219.9998

$ws.Range("H71").Value = 1584.2858
$ws.Range("J71").Value = 1532.6666
$ws.Range("L71").Value = 13793.9994
$ws.Range("N71").Value = -21905.9994

$ws.Range("H122").Value = 695.9375
$ws.Range("I122").Value = 749.8
$ws.Range("J122").Value = 671.4545000000001
$ws.Range("K122").Value = 6748.2
$ws.Range("L122").Value = 6043.0905
$ws.Range("M122").Value = -4298.2
$ws.Range("N122").Value = -10943.0905

$ws.Range("H139").Value = 5124.125
$ws.Range("I139").Value = 3198.6
$ws.Range("K139").Value = 9595.799999999999
$ws.Range("M139").Value = -4455.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 1651.5
$ws.Range("I12").Value = 1651.5
$ws.Range("K12").Value = 1651.5
$ws.Range("M12").Value = -1511.5

$ws.Range("H122").Value = 7355252
$ws.Range("I122").Value = 8930216
$ws.Range("K122").Value = 26790648
$ws.Range("M122").Value = -26788198

$ws.Range("H132").Value = 12
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888

$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 2925
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 2925
$ws.Range("L10").Value = 150
$ws.Range("M10").Value = -2785
$ws.Range("N10").Value = -430

$ws.Range("H12").Value = 883.5714
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 1544
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1544
$ws.Range("M12").Value = 167
$ws.Range("N12").Value = -1884

$ws.Range("H16").Value = 342.66666
$ws.Range("I16").Value = 342.66666
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 342.66666
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -172.66666
$ws.Range("N16").ClearContents()

$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864

$ws.Range("H46").Value = 3446.9565
$ws.Range("I46").Value = 2548.625
$ws.Range("K46").Value = 2548.625
$ws.Range("M46").Value = -2360.625

$ws.Range("H55").Value = 323.07144
$ws.Range("I55").Value = 266.77777
$ws.Range("J55").Value = 424.4
$ws.Range("K55").Value = 266.77777
$ws.Range("L55").Value = 424.4
$ws.Range("M55").Value = -93.77776999999998
$ws.Range("N55").Value = -770.4

$ws.Range("H68").Value = 1820.2
$ws.Range("I68").Value = 1820.2
$ws.Range("K68").Value = 1820.2
$ws.Range("M68").Value = -1071.2

$ws.Range("H71").Value = 1820.2
$ws.Range("I71").Value = 1820.2
$ws.Range("K71").Value = 9101
$ws.Range("M71").Value = -5357

$ws.Range("H122").Value = 3216
$ws.Range("I122").Value = 3168.6667
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 9506.000100000001
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -7056.000100000001
$ws.Range("N122").Value = -15400

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws.Range("H132").Value = 7250.6665
$ws.Range("I132").Value = 2004
$ws.Range("K132").Value = 6012
$ws.Range("M132").Value = -3482

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 250
$ws.Range("J10").Value = 250
$ws.Range("L10").Value = 250
$ws.Range("N10").Value = -588

$ws.Range("H18").Value = 3642.4285
$ws.Range("I18").Value = 5999.25
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 5999.25
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = -5826.25
$ws.Range("N18").Value = -846

$ws.Range("H122").Value = 1814.5
$ws.Range("I122").Value = 1814.5
$ws.Range("K122").Value = 5443.5
$ws.Range("M122").Value = -2993.5

$ws.Range("H126").Value = 1499.5
$ws.Range("I126").Value = 1499.5
$ws.Range("K126").Value = 4498.5
$ws.Range("M126").Value = -2028.5

$ws.Range("H136").Value = 1673.88
$ws.Range("I136").Value = 1221.5238
$ws.Range("K136").Value = 3664.5714
$ws.Range("M136").Value = -1114.5714
